$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.891.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.08%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.811.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.80%  "
$ws.Range("E4").Value = "  +0.55%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.82%  "
$ws.Range("E6").Value = "  +0.49%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4292"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.50%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3701"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07236"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8653"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "2.061.88"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +18.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.32"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.73%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.633"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.95%  "
$ws.Range("E14").Value = "  +3.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06935"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "80.83"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.18%  "
$ws.Range("E17").Value = "  +1.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008912"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.97%  "
$ws.Range("E19").Value = "  +0.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.925.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.194"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.276.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +15.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.74%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.885"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.04%  "
$ws.Range("E27").Value = "  +1.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.244"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.921"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +15.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "114.64"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08955"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7438"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.80%  "
$ws.Range("E33").Value = "  +4.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.432"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.01%  "
$ws.Range("E35").Value = "  +3.19%  "
$ws.Range("E36").Value = "  +0.61%  "
$ws.Range("E37").Value = "  +5.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05233"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01926"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5095"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.748"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +11.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1652"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.484"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.290"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "107.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.39"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.49%  "
$ws.Range("E47").Value = "  +0.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.649"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4560"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06266"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.801"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.08%  "
